$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 31

# Columns A-D hold text values (date, time, weekday, zero-padded week).
# Force text number format so Excel doesn't auto-convert them to
# dates/numbers and strip the leading zero, then restore the default
# "Normal" style (no custom number format) to match the rest of the sheet.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 2).NumberFormat = "@"
$ws.Cells.Item($row, 3).NumberFormat = "@"
$ws.Cells.Item($row, 4).NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = "2025-01-17"
$ws.Cells.Item($row, 2).Value = "22:52:10"
$ws.Cells.Item($row, 3).Value = "Friday"
$ws.Cells.Item($row, 4).Value = "02"

$ws.Cells.Item($row, 1).Style = "Normal"
$ws.Cells.Item($row, 2).Style = "Normal"
$ws.Cells.Item($row, 3).Style = "Normal"
$ws.Cells.Item($row, 4).Style = "Normal"

$ws.Cells.Item($row, 5).Value = 126834
$ws.Cells.Item($row, 6).Value = 142071
$ws.Cells.Item($row, 7).Value = 169319
$ws.Cells.Item($row, 8).Value = 158476
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 143050
$ws.Cells.Item($row, 11).Value = -1
$ws.Cells.Item($row, 12).Value = -1
$ws.Cells.Item($row, 13).Value = 192252
$ws.Cells.Item($row, 14).Value = 115535
$ws.Cells.Item($row, 15).Value = 45519
$ws.Cells.Item($row, 16).Value = 28511
$ws.Cells.Item($row, 17).Value = 65790
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 49533
$ws.Cells.Item($row, 20).Value = -1
